$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Process edits from the bottom of the document upward so that paragraph
# indices for not-yet-processed (earlier) content remain valid.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Change: "Front-end" subtitle - merge split runs ("Front-" + "end") into a
# single run and drop the spell-check proofErr markers.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(46)
$rng = $p.Range
$xml = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="Subttulo"/>
    <w:spacing w:line="240" w:lineRule="auto"/>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:eastAsia="Times New Roman" w:hAnsi="Aptos"/>
      <w:b/>
      <w:color w:val="000000"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:bookmarkStart w:id="3" w:name="_shx0bhwy5wap"/>
  <w:bookmarkEnd w:id="3"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:eastAsia="Times New Roman" w:hAnsi="Aptos"/>
      <w:b/>
      <w:color w:val="000000"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>Front-end</w:t>
  </w:r>
</w:p>
"@
$rng.InsertXML($xml)

# ---------------------------------------------------------------------------
# Change: paragraph right before the blank "ind firstLine=708" run of blank
# paragraphs used to be styled "PargrafodaLista" with ind left=1440; it now
# matches the plain blank paragraphs (ind firstLine=708), and two additional
# blank paragraphs of the same style are inserted right after it.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(41)
$rng = $p.Range
$blankPPr = @"
  <w:pPr>
    <w:ind w:firstLine="708"/>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
"@
$xml = @"
<w:p $wns>
$blankPPr
</w:p>
<w:p $wns>
$blankPPr
</w:p>
<w:p $wns>
$blankPPr
</w:p>
"@
$rng.InsertXML($xml)

# ---------------------------------------------------------------------------
# Change: insert a new blank paragraph (ind left=708) right after the
# "Pontuação total acumulada..." bullet and before the existing blank
# paragraph that follows it.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(26)
$rng = $p.Range
$rng.Collapse(0)
$xml = @"
<w:p $wns>
  <w:pPr>
    <w:ind w:left="708"/>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
"@
$rng.InsertXML($xml)

# ---------------------------------------------------------------------------
# Change: "Porcentagem total acumulada..." -> "Pontuação total acumulada..."
# (only the first three runs that spell "Porcentagem " are merged into one
# run reading "Pontuação "; the remainder of the sentence is untouched).
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(26)
$rng = $p.Range
$xml = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="PargrafodaLista"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="19"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">Pontuação </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">total </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>acumulada</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> com a conclusão das atividades</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>.</w:t>
  </w:r>
</w:p>
"@
$rng.InsertXML($xml)

# ---------------------------------------------------------------------------
# Change: "... através de uma porcentagem a medida que as tarefas são
# concluídas." -> "... através de uma porcentagem de conclusão."
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(25)
$rng = $p.Range
$xml = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="PargrafodaLista"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="19"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Exibição do progresso geral na jornada</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> através de uma porcentagem de conclusão</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>.</w:t>
  </w:r>
</w:p>
"@
$rng.InsertXML($xml)

# ---------------------------------------------------------------------------
# Change: remove the whole "Interface Lúdica de Progresso" bullet section
# (heading + 2 bullets) and leave a single empty, justified paragraph.
# ---------------------------------------------------------------------------
$pStart = $d.Paragraphs.Item(18)
$pEnd = $d.Paragraphs.Item(20)
$rng = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$xml = @"
<w:p $wns>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
"@
$rng.InsertXML($xml)

# ---------------------------------------------------------------------------
# Change: "História de Usuário" story paragraph - drop the justified
# alignment and merge the split runs (caused by a grammar-check proofErr
# marker around "Eu") back into a single run.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(4)
$rng = $p.Range
$xml = @"
<w:p $wns>
  <w:pPr>
    <w:spacing w:after="160" w:line="259" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Como um usuário do sistema em busca de autoaperfeiçoamento e progresso, Eu desejo uma visualização clara e abrangente do meu progresso ao longo da jornada do herói, Para ter uma compreensão detalhada das tarefas concluídas e em andamento, mantendo-me motivado e focado em meus objetivos.</w:t>
  </w:r>
</w:p>
"@
$rng.InsertXML($xml)
